$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5704150199890137
$ws.Range("G3").Value = 0.1509251594543457
$ws.Range("G4").Value = 2.053256988525391
$ws.Range("G5").Value = 0.04685783386230469
$ws.Range("G6").Value = 9.754637956619263
$ws.Range("G7").Value = 0.01299595832824707
$ws.Range("G8").Value = 0.1078310012817383
$ws.Range("G9").Value = 1.276360988616943
$ws.Range("G10").Value = 5.714885950088501
$ws.Range("G11").Value = 0.3073081970214844
$ws.Range("G12").Value = 0.01186394691467285
$ws.Range("G13").Value = 0.03953409194946289
$ws.Range("G14").Value = 14.33312797546387
$ws.Range("G15").Value = 99.24435091018677
$ws.Range("G16").Value = 2.313521146774292
$ws.Range("G17").Value = 0.01349592208862305
$ws.Range("G18").Value = 0.08600306510925293
$ws.Range("G19").Value = 0.4214069843292236
$ws.Range("G20").Value = 0.6422049999237061
$ws.Range("G21").Value = 3.964146852493286
$ws.Range("G22").Value = 0.09005188941955566
$ws.Range("G23").Value = 0.01261591911315918
$ws.Range("G24").Value = 0.03083395957946777
$ws.Range("G25").Value = 0.02793502807617188
$ws.Range("G26").Value = 0.0130620002746582
$ws.Range("G27").Value = 0.06688213348388672
$ws.Range("G28").Value = 0.01235103607177734
$ws.Range("G29").Value = 0.1533589363098145
$ws.Range("G30").Value = 2.0609290599823
$ws.Range("G31").Value = 0.493441104888916
$ws.Range("G32").Value = 0.03567218780517578
$ws.Range("G33").Value = 9.741075038909912
$ws.Range("G34").Value = 0.01178598403930664
$ws.Range("G35").Value = 0.1100990772247314
$ws.Range("G36").Value = 0.04321789741516113
$ws.Range("G37").Value = 1.743595838546753
$ws.Range("G38").Value = 0.4172549247741699
$ws.Range("G39").Value = 25.06800580024719
$ws.Range("G40").Value = 179.1279811859131
$ws.Range("G41").Value = 1.26123309135437
$ws.Range("G42").Value = 0.1114749908447266
$ws.Range("G43").Value = 6.47984504699707
$ws.Range("G44").Value = 0.04830098152160645
$ws.Range("G45").Value = 0.406527042388916
$ws.Range("G46").Value = 0.842940092086792
$ws.Range("G47").Value = 2.003684997558594
$ws.Range("G48").Value = 0.7498798370361328
$ws.Range("G49").Value = 0.03198504447937012
$ws.Range("G50").Value = 0.1816260814666748
$ws.Range("G51").Value = 0.01298999786376953
$ws.Range("G52").Value = 0.05528092384338379
$ws.Range("G53").Value = 0.2879509925842285
$ws.Range("G54").Value = 0.01412296295166016
$ws.Range("G55").Value = 0.07988715171813965
$ws.Range("G56").Value = 0.0333409309387207
$ws.Range("G57").Value = 5.730488061904907
$ws.Range("G58").Value = 0.2275040149688721
$ws.Range("G59").Value = 53.30519104003906
$ws.Range("G60").Value = 0.0175929069519043
$ws.Range("G61").Value = 0.06056118011474609
$ws.Range("G62").Value = 3.676536083221436
$ws.Range("G63").Value = 0.5459098815917969
$ws.Range("G64").Value = 1.279531002044678
$ws.Range("G65").Value = 0.01209402084350586
$ws.Range("G66").Value = 0.04854607582092285
$ws.Range("G67").Value = 0.1074941158294678
$ws.Range("G68").Value = 0.3114171028137207
$ws.Range("G69").Value = 5.765976905822754
$ws.Range("G70").Value = 0.01941585540771484
$ws.Range("G71").Value = 0.3508529663085938
$ws.Range("G72").Value = 1.050707101821899
$ws.Range("G73").Value = 0.08148503303527832
$ws.Range("G74").Value = 0.1732320785522461
